$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": insert a new advisor row at row 28 (all values 0),
# pushing the existing rows (28..45) down to (29..46).
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(28).Insert()
$ws1.Range("A28").Value = "OFICINA-CATAECSA"
$ws1.Range("B28").Value = "LOAIZA TINOCO JUAN PABLO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(28, $col).Value = 0
}

# The trailing "x de 43" summary row (now row 46) counts advisors; update the
# count to reflect the newly-added advisor (43 -> 44).
for ($col = 3; $col -le 18; $col++) {
    $oldText = $ws1.Cells.Item(46, $col).Value()
    $newText = $oldText.Replace("de 43", "de 44")
    $ws1.Cells.Item(46, $col).Value = $newText
}

# Sheet "VENTA MENSUAL": same insert at row 28.
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(28).Insert()
$ws2.Range("A28").Value = "OFICINA-CATAECSA"
$ws2.Range("B28").Value = "LOAIZA TINOCO JUAN PABLO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(28, $col).Value = 0
}
